$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "total change" column (G) values for rows 2-15 from 0.11 to 0.0
$ws.Range("G2:G15").Value = 0.0
